$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-8 for columns T, U (and H/V/W shared-range churn) per the new
#     "try the runs again" parameters (12-month overhaul hours, 100-aircraft fleet
#     size for the 'et_av' column) before dropping the now-unused rows 9:12. ---

# Row 2
$ws.Range("T2").Formula = "=24*30*12"
$ws.Range("U2").Formula = "=100*365*24"

# Row 3
$ws.Range("U3").Formula = "=100*365*24"

# Row 4
$ws.Range("T4").Formula = "=24*30*12"
$ws.Range("U4").Formula = "=100*365*24"

# Row 5
$ws.Range("S5").Value = 14400
$ws.Range("T5").Formula = "=24*30*9"
$ws.Range("U5").Formula = "=100*365*24"

# Row 6
$ws.Range("T6").Formula = "=24*30*12"
$ws.Range("U6").Formula = "=100*365*24"

# Row 7
$ws.Range("U7").Formula = "=100*365*24"
$ws.Range("W7").Formula = "=24*365*12"

# Row 8
$ws.Range("S8").Value = 18000
$ws.Range("T8").Formula = "=24*30*12"
$ws.Range("U8").Formula = "=100*365*24"
$ws.Range("W8").Formula = "=24*365*12"

# --- Drop the trailing runs (rows 9-12); only rows 2-8 of data remain ---
$ws.Range("A9:A12").EntireRow.Delete()

# --- Shrink the AutoFilter down to the new data extent ---
$ws.AutoFilterMode = $false
$ws.Range("U1:U13").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$U`$1:`$U`$13"
    }
}

# --- Restore cursor position to where the author last left off ---
$ws.Range("S17").Select() | Out-Null
